$d = $word.ActiveDocument

# 1. Replace the ID placeholder text, absorbing the trailing space run
#    (the space-only run has identical formatting and gets merged/removed
#    when the combined text is replaced without a trailing space).
$d.Content.Find.Execute(
    "**ID__AFFARS_5338_topic_1__ID** ", $true, $false, $false, $false,
    $false, $true, 1, $false, "**ID__AFFARS_PART_5338__ID**", 2
)

# 2. Update the first paragraph's formatting: add a paragraph border
#    (5pt space on every side) and change the left indent from 120 to 225
#    twips (i.e. 6pt -> 11.25pt).
$p = $d.Paragraphs(1)
$p.Borders.DistanceFromTop = 5
$p.Borders.DistanceFromLeft = 5
$p.Borders.DistanceFromBottom = 5
$p.Borders.DistanceFromRight = 5
$p.Range.ParagraphFormat.LeftIndent = 11.25
